$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 358, shifting existing rows 358..382 down to 359..383
$ws.Rows.Item(358).Insert()

$r = 358
$ws.Cells.Item($r, 1).Value = 9
$ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($r, 3).Value = "Metropolitana"
$ws.Cells.Item($r, 4).Value = 44714
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r + 1, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = 13
$ws.Cells.Item($r, 6).Value = 100112044
$ws.Cells.Item($r, 7).Value = "Perejil"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 79
$ws.Cells.Item($r, 11).Value = 8000
$ws.Cells.Item($r, 12).Value = 9000
$ws.Cells.Item($r, 13).Value = 8494
$ws.Cells.Item($r, 14).Value = "$/docena de atados"
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 2831
$ws.Cells.Item($r, 17).Value = 3
$ws.Cells.Item($r, 18).Value = "Hortaliza"
